$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.737.15"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "1.846.65"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.41"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4298"
$ws.Range("E7").Value = "  +1.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3660"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.08"
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07343"
$ws.Range("E10").Value = "  +0.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8770"
$ws.Range("E11").Value = "  -2.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.75"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "1.812.34"
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.335"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.520"
$ws.Range("E15").Value = "  -0.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06961"
$ws.Range("E16").Value = "  +1.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "79.82"
$ws.Range("E18").Value = "  +1.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009011"
$ws.Range("E19").Value = "  +1.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9998"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("D22").Value = "27.539.54"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.972"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.37"
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("D25").Value = "1.998.23"
$ws.Range("E25").Value = "  -3.24%  "
$ws.Range("E26").Value = "  -3.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.33"
$ws.Range("E27").Value = "  +1.37%  "
$ws.Range("E28").Value = "  +2.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "119.83"
$ws.Range("E29").Value = "  +7.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.257"
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.866"
$ws.Range("E31").Value = "  +2.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08906"
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7523"
$ws.Range("E33").Value = "  -2.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.552"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.972"
$ws.Range("E35").Value = "  +1.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.125"
$ws.Range("E36").Value = "  +3.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.106"
$ws.Range("E37").Value = "  +0.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05424"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01933"
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.855"
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5086"
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1656"
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.636"
$ws.Range("E43").Value = "  -2.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.331"
$ws.Range("E44").Value = "  +1.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06542"
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.28"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "104.88"
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4663"
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.626"
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.16"
$ws.Range("E51").Value = "  -0.09%  "
